$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BB mirrors column BA (same style/format), then gets its own
# values: the header date moves one quarter forward, most rows simply repeat
# the BA forecast value, and the last three rows (19-21) get the newly
# evaluated EQUIPMENT forecast numbers.

# Copy BA's formatting for the header and the data rows that receive a value
# in column BB (rows 3-21), so newly written cells inherit the correct
# number formats/styles without creating new style entries. Rows 2 and 22
# stay untouched (no BA/BB value there).
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)

$ws.Range("BA3:BA21").Copy()
$ws.Range("BB3:BB21").PasteSpecial(-4122)

# Header date for the new column.
$ws.Range("BB1").Value = 45986

# Column BB body values: mostly a copy of column BA, except rows 19-21 which
# carry updated forecast values (new EQUIPMENT evaluation data).
$bbValues = @{
    3  = -5.109987415979145
    4  = 2.253603114136604
    5  = 3.738382206110891
    6  = 1.165974434765671
    7  = -0.0426719751787874
    8  = 1.529758493743438
    9  = 1.358758534900462
    10 = 1.664905435092301
    11 = 2.145670176886982
    12 = 1.976124254426503
    13 = 0.7060158009350337
    14 = -4.180878843351332
    15 = 1.312484974417294
    16 = 2.386394320099283
    17 = 0.2104414886460626
    18 = -0.3095793941792935
    19 = -0.08656168856399082
    20 = 0.6232357314897463
    21 = 0.7596754881313172
}

foreach ($row in $bbValues.Keys) {
    $ws.Cells.Item($row, 54).Value = $bbValues[$row]
}
